$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "School Name"
$ws.Range("C2").Value = "Dr GR PUBLIC SCHOOL"

$ws.Columns.Item(3).ColumnWidth = 21.5
